# Edit applied per commit: "Default Log Path Fixed."
# - Rename column A header from "Time(sec)" to "Time"
# - Update logged sensor data (time/pressure columns, RPM columns) for existing rows
#   and append 6 additional rows of newly logged data (rows 13-18)
# - Adjust the saved workbook window position/size metadata

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "Time(sec)" -> "Time" ---
$ws.Range("A1").Value = "Time"

# --- Updated data rows (A: time, B: breath, C: pressure, D: LRPM, E: RRPM) ---
$ws.Range("A2").Value = 12.327
$ws.Range("C2").Value = 10998.64453125
$ws.Range("A3").Value = 12.586
$ws.Range("C3").Value = 10998.6728515625
$ws.Range("A4").Value = 12.844
$ws.Range("C4").Value = 10998.5966796875
$ws.Range("D4").Value = 120
$ws.Range("E4").Value = 120
$ws.Range("A5").Value = 13.103
$ws.Range("C5").Value = 10998.7626953125
$ws.Range("D5").Value = 120
$ws.Range("E5").Value = 120
$ws.Range("A6").Value = 13.361
$ws.Range("C6").Value = 10998.654296875
$ws.Range("D6").Value = 120
$ws.Range("E6").Value = 120
$ws.Range("A7").Value = 13.62
$ws.Range("C7").Value = 10998.646484375
$ws.Range("D7").Value = 120
$ws.Range("E7").Value = 120
$ws.Range("A8").Value = 13.879
$ws.Range("C8").Value = 10998.6025390625
$ws.Range("A9").Value = 14.137
$ws.Range("C9").Value = 10998.6123046875
$ws.Range("A10").Value = 14.396
$ws.Range("C10").Value = 10998.640625
$ws.Range("A11").Value = 14.654
$ws.Range("C11").Value = 10998.673828125
$ws.Range("A12").Value = 14.913
$ws.Range("C12").Value = 10998.5576171875
$ws.Range("A13").Value = 15.171
$ws.Range("B13").Value = -0.016599999740719795
$ws.Range("C13").Value = 10998.5693359375
$ws.Range("D13").Value = 120
$ws.Range("E13").Value = 120
$ws.Range("A14").Value = 15.43
$ws.Range("B14").Value = -0.016599999740719795
$ws.Range("C14").Value = 10998.640625
$ws.Range("D14").Value = 120
$ws.Range("E14").Value = 120
$ws.Range("A15").Value = 15.689
$ws.Range("B15").Value = -0.016599999740719795
$ws.Range("C15").Value = 10998.6142578125
$ws.Range("D15").Value = 120
$ws.Range("E15").Value = 120
$ws.Range("A16").Value = 15.947
$ws.Range("B16").Value = -0.016599999740719795
$ws.Range("C16").Value = 10998.576171875
$ws.Range("D16").Value = 120
$ws.Range("E16").Value = 120
$ws.Range("A17").Value = 16.206
$ws.Range("B17").Value = -0.016599999740719795
$ws.Range("C17").Value = 10998.5732421875
$ws.Range("D17").Value = 120
$ws.Range("E17").Value = 120
$ws.Range("A18").Value = 16.464
$ws.Range("B18").Value = -0.016599999740719795
$ws.Range("C18").Value = 10998.5048828125
$ws.Range("D18").Value = 120
$ws.Range("E18").Value = 120

# --- Window view metadata (saved window position / size) ---
$win = $wb.Windows.Item(1)
$win.Left = 480
$win.Top = 84
$win.Width = 22056
$win.Height = 11424
